$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)
$ws2.Name = "工作表1"
$ws2.Range("B1").Value = "Quick sort"

$co = $ws2.ChartObjects().Add(100, 100, 400, 300)
$chart = $co.Chart
$chart.ChartType = 75
$s1 = $chart.SeriesCollection().NewSeries()

$tests = @("=工作表1!`$B`$1", '=工作表1!$B$1', "=Sheet1!`$B`$1")
foreach ($t in $tests) {
  Write-Host "Trying: $t"
  try {
    $s1.Name = $t
    Write-Host "  -> ok"
  } catch {
    Write-Host "  -> error: $_"
  }
}
